# ============================================================================
# Add 2022-Q4 data:
#  - Insert a new "2022-Q4" worksheet right after "总计" (总计 stays sheet 1,
#    existing "2022-Q3"/"2022-Q1"/"2021-Q4" sheets shift right by one).
#  - Populate the new sheet with the Q4 fund-holding table.
#  - Prepend a "2022-Q4" summary row to the "总计" sheet and renumber the
#    existing index column (A) for the rows that shift down.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1) "总计" (summary) sheet: add the 2022-Q4 row at the top of the table,
#    push 2022-Q3 / 2022-Q1 / 2021-Q4 down one row each.
# ----------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Copy the existing index-column style (bold/centered) onto the new A5 cell
# before we populate it, since it currently has no explicit style.
$summary.Range("A2").Copy()
$summary.Range("A5").PasteSpecial(-4122)

# Write bottom-up so every row ends on its final value in one shot.
$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 4
$summary.Range("D5").Value = 3.61

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 3
$summary.Range("D4").Value = 2.59

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 20
$summary.Range("D3").Value = 12.07

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 22
$summary.Range("D2").Value = 13.2


# ----------------------------------------------------------------------
# 2) Insert the new "2022-Q4" worksheet right after "总计".
# ----------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q4"
$ws = $newSheet

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Header row + index column (A) share the bold/centered style used on the
# other quarterly sheets; copy it in from the "总计" sheet's header cell.
$summary.Range("B1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

# Force the numeric-looking text columns (fund code + the four formatted
# decimal columns) to be stored as text, matching the source data, instead
# of being auto-coerced to numbers by the COM value-assignment heuristics.
$codeFmt = $ws.Range("B2:B23")
$decFmt = $ws.Range("D2:G23")
$codeFmt.NumberFormat = "@"
$decFmt.NumberFormat = "@"


# ----------------------------------------------------------------------
# 3) Data rows (22 funds), A1:H23 overall.
# ----------------------------------------------------------------------
$ws.Range("B2").Value = "001511"
$ws.Range("C2").Value = "兴全新视野灵活配置定期开放混合"
$ws.Range("D2").Value = "111.90"
$ws.Range("E2").Value = "90.77"
$ws.Range("F2").Value = "3.57"
$ws.Range("G2").Value = "3.9948"
$ws.Range("H2").Value = 6
$ws.Range("B3").Value = "163415"
$ws.Range("C3").Value = "兴全商业模式优选混合（LOF）"
$ws.Range("D3").Value = "108.38"
$ws.Range("E3").Value = "93.15"
$ws.Range("F3").Value = "3.50"
$ws.Range("G3").Value = "3.7933"
$ws.Range("H3").Value = 6
$ws.Range("B4").Value = "011056"
$ws.Range("C4").Value = "博时汇兴回报一年持有期灵活配置混合"
$ws.Range("D4").Value = "91.53"
$ws.Range("E4").Value = "68.57"
$ws.Range("F4").Value = "2.92"
$ws.Range("G4").Value = "2.6727"
$ws.Range("H4").Value = 7
$ws.Range("B5").Value = "013797"
$ws.Range("C5").Value = "博时优质鑫选一年持有期混合A"
$ws.Range("D5").Value = "38.62"
$ws.Range("E5").Value = "79.82"
$ws.Range("F5").Value = "3.04"
$ws.Range("G5").Value = "1.1740"
$ws.Range("H5").Value = 6
$ws.Range("B6").Value = "011177"
$ws.Range("C6").Value = "博时汇融回报一年持有期混合A"
$ws.Range("D6").Value = "24.28"
$ws.Range("E6").Value = "77.15"
$ws.Range("F6").Value = "2.35"
$ws.Range("G6").Value = "0.5706"
$ws.Range("H6").Value = 9
$ws.Range("B7").Value = "001236"
$ws.Range("C7").Value = "博时丝路主题股票A"
$ws.Range("D7").Value = "10.52"
$ws.Range("E7").Value = "87.16"
$ws.Range("F7").Value = "3.14"
$ws.Range("G7").Value = "0.3303"
$ws.Range("H7").Value = 7
$ws.Range("B8").Value = "009740"
$ws.Range("C8").Value = "博时研究臻选三年持有期灵活配置混合A"
$ws.Range("D8").Value = "7.01"
$ws.Range("E8").Value = "80.20"
$ws.Range("F8").Value = "3.23"
$ws.Range("G8").Value = "0.2264"
$ws.Range("H8").Value = 6
$ws.Range("B9").Value = "001468"
$ws.Range("C9").Value = "广发改革先锋灵活配置混合"
$ws.Range("D9").Value = "5.72"
$ws.Range("E9").Value = "93.29"
$ws.Range("F9").Value = "2.94"
$ws.Range("G9").Value = "0.1682"
$ws.Range("H9").Value = 5
$ws.Range("B10").Value = "011845"
$ws.Range("C10").Value = "博时周期优选混合A"
$ws.Range("D10").Value = "1.97"
$ws.Range("E10").Value = "81.65"
$ws.Range("F10").Value = "3.28"
$ws.Range("G10").Value = "0.0646"
$ws.Range("H10").Value = 6
$ws.Range("B11").Value = "015031"
$ws.Range("C11").Value = "博时远见回报混合C"
$ws.Range("D11").Value = "1.05"
$ws.Range("E11").Value = "81.10"
$ws.Range("F11").Value = "3.20"
$ws.Range("G11").Value = "0.0336"
$ws.Range("H11").Value = 6
$ws.Range("B12").Value = "011340"
$ws.Range("C12").Value = "博时战略新材料主题混合A"
$ws.Range("D12").Value = "0.84"
$ws.Range("E12").Value = "79.48"
$ws.Range("F12").Value = "3.76"
$ws.Range("G12").Value = "0.0316"
$ws.Range("H12").Value = 4
$ws.Range("B13").Value = "002556"
$ws.Range("C13").Value = "博时丝路主题股票C"
$ws.Range("D13").Value = "0.92"
$ws.Range("E13").Value = "87.16"
$ws.Range("F13").Value = "3.14"
$ws.Range("G13").Value = "0.0289"
$ws.Range("H13").Value = 7
$ws.Range("B14").Value = "015030"
$ws.Range("C14").Value = "博时远见回报混合A"
$ws.Range("D14").Value = "0.69"
$ws.Range("E14").Value = "81.10"
$ws.Range("F14").Value = "3.20"
$ws.Range("G14").Value = "0.0221"
$ws.Range("H14").Value = 6
$ws.Range("B15").Value = "011341"
$ws.Range("C15").Value = "博时战略新材料主题混合C"
$ws.Range("D15").Value = "0.56"
$ws.Range("E15").Value = "79.48"
$ws.Range("F15").Value = "3.76"
$ws.Range("G15").Value = "0.0211"
$ws.Range("H15").Value = 4
$ws.Range("B16").Value = "014212"
$ws.Range("C16").Value = "博时研究优享混合A"
$ws.Range("D16").Value = "0.59"
$ws.Range("E16").Value = "79.20"
$ws.Range("F16").Value = "3.19"
$ws.Range("G16").Value = "0.0188"
$ws.Range("H16").Value = 7
$ws.Range("B17").Value = "009741"
$ws.Range("C17").Value = "博时研究臻选三年持有期灵活配置混合C"
$ws.Range("D17").Value = "0.45"
$ws.Range("E17").Value = "80.20"
$ws.Range("F17").Value = "3.23"
$ws.Range("G17").Value = "0.0145"
$ws.Range("H17").Value = 6
$ws.Range("B18").Value = "014913"
$ws.Range("C18").Value = "博时研究回报混合A"
$ws.Range("D18").Value = "0.47"
$ws.Range("E18").Value = "77.91"
$ws.Range("F18").Value = "2.63"
$ws.Range("G18").Value = "0.0124"
$ws.Range("H18").Value = 10
$ws.Range("B19").Value = "013798"
$ws.Range("C19").Value = "博时优质鑫选一年持有期混合C"
$ws.Range("D19").Value = "0.37"
$ws.Range("E19").Value = "79.82"
$ws.Range("F19").Value = "3.04"
$ws.Range("G19").Value = "0.0112"
$ws.Range("H19").Value = 6
$ws.Range("B20").Value = "014914"
$ws.Range("C20").Value = "博时研究回报混合C"
$ws.Range("D20").Value = "0.14"
$ws.Range("E20").Value = "77.91"
$ws.Range("F20").Value = "2.63"
$ws.Range("G20").Value = "0.0037"
$ws.Range("H20").Value = 10
$ws.Range("B21").Value = "011846"
$ws.Range("C21").Value = "博时周期优选混合C"
$ws.Range("D21").Value = "0.09"
$ws.Range("E21").Value = "81.65"
$ws.Range("F21").Value = "3.28"
$ws.Range("G21").Value = "0.0030"
$ws.Range("H21").Value = 6
$ws.Range("B22").Value = "011178"
$ws.Range("C22").Value = "博时汇融回报一年持有期混合C"
$ws.Range("D22").Value = "0.12"
$ws.Range("E22").Value = "77.15"
$ws.Range("F22").Value = "2.35"
$ws.Range("G22").Value = "0.0028"
$ws.Range("H22").Value = 9
$ws.Range("B23").Value = "014213"
$ws.Range("C23").Value = "博时研究优享混合C"
$ws.Range("D23").Value = "0.08"
$ws.Range("E23").Value = "79.20"
$ws.Range("F23").Value = "3.19"
$ws.Range("G23").Value = "0.0026"
$ws.Range("H23").Value = 7

# Text-format scratch ranges are no longer needed once the values are in;
# clear them so the cells fall back to the workbook's default (un-styled)
# formatting, matching the source file.
$codeFmt.ClearFormats()
$decFmt.ClearFormats()

# Index column A2:A23 uses the same bold/centered style as the header and
# the "总计" index column.
$summary.Range("A2").Copy()
$ws.Range("A2:A23").PasteSpecial(-4122)
for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

